# Add YouTube links to help and user guide for collaboration features.
#
# The three screenshots below were (re)inserted/touched as part of adding
# the YouTube demo links, so Word now marks their runs as "do not spell/
# grammar check" (<w:noProof/>) -- the same way it already marks the
# existing YouTube-link image and the "Collaborate" screenshot further
# down in the document. We reproduce that by setting NoProofing = True on
# each picture's Range, which is the InlineShape/Range property backing
# the <w:rPr><w:noProof/></w:rPr> markup.

$d = $word.ActiveDocument

# Picture 2: screenshot after "Submit/Share Package" (anchorId 4885995B)
$d.InlineShapes.Item(2).Range.NoProofing = $true

# Picture 3: screenshot after "Submit Package to EDI" (anchorId 38388EE5)
$d.InlineShapes.Item(3).Range.NoProofing = $true

# Picture 5: screenshot of the package-errors warning page (anchorId 1C12CD13)
$d.InlineShapes.Item(5).Range.NoProofing = $true
